$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction in SA algorithm log values (run_15): collapse fitness values
# in column C to flat plateaus for generations 0-91 (rows 2-93).
$ws.Range("C2:C31").Value = 7917
$ws.Range("C32:C78").Value = 7318
$ws.Range("C79:C93").Value = 7293
